# Update the dSF column (F) values for the rows that were repulled from
# the source data. All other columns/rows are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -3
    11 = -2
    12 = -2
    15 = -1
    17 = -3
    18 = -1
    28 = -4
    33 = -1
    34 = 4
    35 = 4
    38 = 0
    39 = 0
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
